# Fix nested array error:
#  - Rename the "data" root prefix to "myData_isArray" throughout the header row
#  - Add new header columns (H:L) describing a nested array ("array_ja_isArray")
#  - Insert a new row 2 holding the array length marker "[4]" for the top-level array,
#    shifting the previous data rows (dog/bird/cat) down by one row
#  - Add the nested-array detail rows (J6:K7) and trailing "ok"/"[2]"/age markers on row 5
#  - Append a new data row (bird2) at the bottom
#  - Update the sheet dimension / ignoredErrors sqref to cover the new A1:L9 used range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): rename "data" prefix -> "myData_isArray", add new columns H:L ---
$ws.Range("A1").Value = "myData_isArray"
$ws.Range("B1").Value = "myData_isArray.name"
$ws.Range("C1").Value = "myData_isArray.breed"
$ws.Range("D1").Value = "myData_isArray.age"
$ws.Range("E1").Value = "myData_isArray.origin"
$ws.Range("F1").Value = "myData_isArray.origin.country"
$ws.Range("G1").Value = "myData_isArray.origin.city"
$ws.Range("H1").Value = "myData_isArray.test"
$ws.Range("I1").Value = "myData_isArray.array_ja_isArray"
$ws.Range("J1").Value = "myData_isArray.array_ja_isArray.test"
$ws.Range("K1").Value = "myData_isArray.array_ja_isArray.test2"
$ws.Range("L1").Value = "myData_isArray.test_age"

# --- Row 2 becomes the array-length marker; the old row-2 animal data (dog1) moves to row 3 ---
$ws.Range("B2:G2").ClearContents()
$ws.Range("A2").Value = "[4]"

# --- Row 3: dog1 (previously on row 2) ---
$ws.Range("B3").Value = "dog1"
$ws.Range("C3").Value = "dog"
$ws.Range("D3").Value = 2
$ws.Range("F3").Value = "TH"
$ws.Range("G3").Value = "BKK"

# --- Row 4: bird1 (previously on row 3); origin.country no longer populated ---
$ws.Range("B4").Value = "bird1"
$ws.Range("C4").Value = "bird"
$ws.Range("D4").Value = 1
$ws.Range("F4").ClearContents()
$ws.Range("G4").Value = "BKK"

# --- Row 5: cat1 (previously on row 4); origin.city corrected to BKK; new nested-array fields ---
$ws.Range("B5").Value = "cat1"
$ws.Range("C5").Value = "cat"
$ws.Range("D5").Value = 7
$ws.Range("F5").Value = "TH"
$ws.Range("G5").Value = "BKK"
$ws.Range("H5").Value = "ok"
$ws.Range("I5").Value = "[2]"
$ws.Range("L5").Value = 1232

# --- Rows 6 & 7: nested array_ja_isArray detail rows ---
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = "qwerty"
$ws.Range("J7").Value = 2
$ws.Range("K7").Value = "qwerty2"

# --- Row 8: new bird2 data row appended ---
$ws.Range("B8").Value = "bird2"
$ws.Range("C8").Value = "bird"
$ws.Range("D8").Value = 1
$ws.Range("F8").Value = "TH"
